# Update status text "Ready for handoff" -> "In Translation" across all
# worksheets that reference it, and shrink the now-narrower "Status"
# columns to match the autofit width that Excel computes for the shorter
# text.

$wb = $excel.ActiveWorkbook

# Target stored width from the archive report is 13.4101845877511 chars.
# Excel's ColumnWidth setter quantizes to whole-pixel boundaries, so the
# nearest value it can actually produce is reached by requesting 12.5.
$oldText = "Ready for handoff"
$newText = "In Translation"
$newWidth = 12.5

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rowCount = $used.Rows.Count
    $colCount = $used.Columns.Count
    for ($r = 1; $r -le $rowCount; $r++) {
        for ($c = 1; $c -le $colCount; $c++) {
            $cell = $ws.Cells.Item($r, $c)
            $text = [string]$cell.Text
            if ($text -eq $oldText) {
                $cell.Value = $newText
                $cell.EntireColumn.ColumnWidth = $newWidth
            }
        }
    }
}
